$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) date values from 2024-12-05 (45631) to 2024-12-07 (45633)
# for all data rows (2 through 34).
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 3).Value = 45633
}

# Remove the custom column width/outline definition that previously applied
# to columns V:Z (22-26), restoring them to the default width.
$ws.Range("V1:Z1").EntireColumn.ClearOutline()
$colRange = $ws.Range($ws.Cells.Item(1, 22), $ws.Cells.Item(1, 26)).EntireColumn
$colRange.ColumnWidth = $ws.StandardWidth
$colRange.OutlineLevel = 0
